$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B28 ---
$b28 = @'
 openActionBarOverflowOrOptionsMenu(getInstrumentation().getTargetContext());
'@
$ws.Range("B28").Value = $b28
$b28r2 = $ws.Range("B28").Characters(2, 76)
$b28r2.Font.Name = "Calibri"
$b28r2.Font.Size = 11
$b28r2.Font.Bold = $true

# --- D28 ---
$d28 = @'
solo.sendKey(solo.MENU);

'@
$ws.Range("D28").Value = $d28
$ws.Range("D28").WrapText = $true

# --- B29 ---
$b29 = @'

        ViewInteraction textView2 = onView(
                allOf(withId(android.R.id.testo), withText("TestoDaSelezionare"),
                        childAtPosition(
                                childAtPosition(
 withClassName(is("com.android.internal.view.menu.ListMenuItemView")),
                                        PosizioneElemento),
                                Posizione),
                        isDisplayed()));
        textView2.perform(click());
'@
$ws.Range("B29").Value = $b29
$b29r2 = $ws.Range("B29").Characters(438, 31)
$b29r2.Font.Name = "Calibri"
$b29r2.Font.Size = 11
$b29r2.Font.Bold = $true

# --- D29 ---
$d29 = @'
 
  solo.clickInList(PosizioneNelMenuDaSelezionare, 0);
'@
$ws.Range("D29").Value = $d29
$ws.Range("D29").WrapText = $true

# --- selection / view ---
$ws.Range("D29").Select()
